$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (20) with the next forecast vector entry.
$ws.Range("A20").Value = 45986
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 2.043309689777173
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = 0.9040423720836799

# Match the date formatting/style used by the rest of column A (e.g. A19).
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = 45986
